$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '61.955.69'
Set-TextCell 2 5 '  -0.44%  '

Set-TextCell 3 4 '2.418.08'
Set-TextCell 3 5 '  -0.28%  '

Set-TextCell 4 5 '  +0.05%  '

Set-TextCell 5 4 '562.32'
Set-TextCell 5 5 '  +0.47%  '

Set-TextCell 6 4 '143.04'
Set-TextCell 6 5 '  -0.90%  '

Set-TextCell 7 5 '  +0.06%  '

Set-TextCell 8 4 '0.531'
Set-TextCell 8 5 '  -0.42%  '

Set-TextCell 9 5 '  -0.50%  '

Set-TextCell 10 5 '  -0.87%  '

Set-TextCell 11 4 '5.19'
Set-TextCell 11 5 '  -4.09%  '

Set-TextCell 12 4 '0.349'
Set-TextCell 12 5 '  -1.24%  '

Set-TextCell 13 4 '25.97'
Set-TextCell 13 5 '  -1.10%  '

Set-TextCell 14 5 '  -1.89%  '

Set-TextCell 15 4 '2.857.85'
Set-TextCell 15 5 '  +0.16%  '

Set-TextCell 16 4 '61.861.36'
Set-TextCell 16 5 '  -0.35%  '

Set-TextCell 17 4 '2.405.28'
Set-TextCell 17 5 '  -0.43%  '

Set-TextCell 18 4 '11.36'
Set-TextCell 18 5 '  +1.22%  '

Set-TextCell 19 4 '323.35'
Set-TextCell 19 5 '  -0.46%  '

Set-TextCell 20 4 '6.83'
Set-TextCell 20 5 '  +0.39%  '

Set-TextCell 21 4 '4.13'
Set-TextCell 21 5 '  -1.83%  '

Set-TextCell 22 5 '  -0.01%  '

Set-TextCell 23 4 '66.73'
Set-TextCell 23 5 '  +1.79%  '

Set-TextCell 24 5 '  -0.37%  '

Set-TextCell 25 4 '8.71'
Set-TextCell 25 5 '  -3.16%  '

Set-TextCell 26 4 '552.41'
Set-TextCell 26 5 '  -7.10%  '

Set-TextCell 27 4 '2.536.86'
Set-TextCell 27 5 '  -0.26%  '

Set-TextCell 28 5 '  +0.10%  '

Set-TextCell 29 4 '0.0₃0931'
Set-TextCell 29 5 '  -1.38%  '

Set-TextCell 30 4 '8.18'
Set-TextCell 30 5 '  -1.68%  '

Set-TextCell 31 4 '1.38'
Set-TextCell 31 5 '  -4.55%  '

Set-TextCell 32 5 '  -2.37%  '

Set-TextCell 33 5 '  -1.40%  '

Set-TextCell 34 4 '1.50'
Set-TextCell 34 5 '  -4.16%  '

Set-TextCell 35 4 '1.00'
Set-TextCell 35 5 '  +0.02%  '

Set-TextCell 36 4 '4.73'
Set-TextCell 36 5 '  -1.57%  '

Set-TextCell 37 4 '0.378'
Set-TextCell 37 5 '  -1.77%  '

Set-TextCell 38 4 '153.50'
Set-TextCell 38 5 '  +0.12%  '

Set-TextCell 39 4 '5.41'
Set-TextCell 39 5 '  -5.73%  '

Set-TextCell 40 4 '18.51'
Set-TextCell 40 5 '  -1.26%  '

Set-TextCell 41 4 '1.80'
Set-TextCell 41 5 '  -0.93%  '

Set-TextCell 42 4 '0.999'
Set-TextCell 42 5 '  -0.06%  '

Set-TextCell 43 4 '146.71'
Set-TextCell 43 5 '  -3.10%  '

Set-TextCell 44 4 '2.23'
Set-TextCell 44 5 '  -6.50%  '

Set-TextCell 45 4 '3.63'
Set-TextCell 45 5 '  -0.87%  '

Set-TextCell 46 4 '0.0527'
Set-TextCell 46 5 '  -2.71%  '

Set-TextCell 47 2 'Mantle'
Set-TextCell 47 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 47 4 '0.591'
Set-TextCell 47 5 '  -0.23%  '

Set-TextCell 48 2 'InjectiveProtocol'
Set-TextCell 48 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 48 4 '19.75'
Set-TextCell 48 5 '  -3.15%  '

Set-TextCell 49 4 '0.0920'
Set-TextCell 49 5 '  -0.22%  '

Set-TextCell 50 5 '  -1.36%  '

Set-TextCell 51 2 'BitgetToken'
Set-TextCell 51 3 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
Set-TextCell 51 4 '1.07'
Set-TextCell 51 5 '  +4.33%  '
